$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.093087315559387
$ws.Range("B1").Value = 1.891440391540527
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.242296457290649
$ws.Range("E1").Value = 1.235905170440674
